$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.232.19"
$ws.Range("E2").Value = "  +3.63%  "
$ws.Range("D3").Value = "1.811.78"
$ws.Range("E3").Value = "  +4.80%  "
$ws.Range("E4").Value = "  -0.37%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "329.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.28%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4458"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.68%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3706"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.87%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.78"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07710"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.40%  "
$ws.Range("E11").Value = "  +2.16%  "
$ws.Range("E12").Value = "  -0.30%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.09"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.29%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.301"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.33%  "
$ws.Range("E15").Value = "  +6.67%  "
$ws.Range("D16").Value = "1.846.28"
$ws.Range("E16").Value = "  +6.68%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.04"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +7.81%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001083"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.87%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06551"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +10.43%  "
$ws.Range("E20").Value = "  -0.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.52"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.91%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.227"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.96%  "
$ws.Range("D23").Value = "28.296.77"
$ws.Range("E23").Value = "  +3.66%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.69"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.63%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.138"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -10.59%  "
$ws.Range("E26").Value = "  +3.91%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "155.98"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.40%  "
$ws.Range("D28").Value = "2.037.74"
$ws.Range("E28").Value = "  +5.78%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.317"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.77%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "128.35"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.197"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.903"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09261"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.15%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.656"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.83%  "
$ws.Range("E35").Value = "  +3.96%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02357"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.46%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2185"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.66%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.180"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.10%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06234"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.24%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6577"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.97%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.199"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.45%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.146"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.60%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.000"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.16%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.410"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("E45").Value = "  +3.99%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6091"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.26%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.769"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.18%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "126.95"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.41%  "
$ws.Range("E49").Value = "  +5.26%  "
$ws.Range("E50").Value = "  +6.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06984"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.58%  "
